$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.27"
$ws.Range("E2").Value = "'-0.08%"
$ws.Range("D3").Value = "'32.29"
$ws.Range("E3").Value = "'1.56%"
$ws.Range("D4").Value = "'5.022"
$ws.Range("E4").Value = "'-1.43%"
$ws.Range("D5").Value = "'0.07630"
$ws.Range("E5").Value = "'-2.27%"
$ws.Range("D6").Value = "'1.948"
$ws.Range("E6").Value = "'-13.37%"
$ws.Range("D7").Value = "'7.870"
$ws.Range("E7").Value = "'1.03%"
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = "'0.9176"
$ws.Range("E8").Value = "'0.05%"
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").Value = "'0.1760"
$ws.Range("E9").Value = "'-0.03%"
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = "'0.07837"
$ws.Range("E10").Value = "'3.97%"
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = "'0.08546"
$ws.Range("E11").Value = "'-4.86%"
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = "'0.03156"
$ws.Range("E12").Value = "'4.10%"
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = "'0.09991"
$ws.Range("E13").Value = "'-0.36%"
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = "'0.001525"
$ws.Range("E14").Value = "'1.39%"
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = "'0.005827"
$ws.Range("E15").Value = "'-1.07%"
$ws.Range("B16").Value = 'UpBots'
$ws.Range("C16").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D16").Value = "'0.007498"
$ws.Range("E16").Value = "'2,116.77%"
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = "'3.461"
$ws.Range("E17").Value = "'-0.23%"
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").Value = "'3.780"
$ws.Range("E18").Value = "'-0.98%"
$ws.Range("E19").Value = "'-4.42%"
$ws.Range("D20").Value = "'0.3345"
$ws.Range("E20").Value = "'1.62%"
$ws.Range("D21").Value = "'0.1298"
$ws.Range("E21").Value = "'-2.87%"
$ws.Range("D22").Value = "'4.264"
$ws.Range("E22").Value = "'0.94%"
$ws.Range("D23").Value = "'0.1990"
$ws.Range("E23").Value = "'9.59%"
$ws.Range("D24").Value = "'0.04504"
$ws.Range("E24").Value = "'-1.98%"
$ws.Range("D25").Value = "'0.001222"
$ws.Range("E25").Value = "'-2.19%"
$ws.Range("E26").Value = "'-1.72%"
$ws.Range("D27").Value = "'0.0001251"
$ws.Range("E27").Value = "'0.19%"
$ws.Range("D39").Value = "'0.01706"
$ws.Range("E39").Value = "'-3.93%"
$ws.Range("E40").Value = "'-2.36%"
$ws.Range("D41").Value = "'0.007452"
$ws.Range("E41").Value = "'1.00%"
$ws.Range("E42").Value = "'-0.82%"
$ws.Range("D43").Value = "'0.002331"
$ws.Range("E43").Value = "'6.60%"
$ws.Range("D44").Value = "'0.01048"
$ws.Range("E44").Value = "'2.07%"
$ws.Range("D45").Value = "'0.00006237"
$ws.Range("E45").Value = "'-0.66%"
$ws.Range("E46").Value = "'0.18%"
$ws.Range("D48").Value = "'0.8204"
$ws.Range("E48").Value = "'11.32%"
$ws.Range("E49").Value = "'0.18%"
$ws.Range("E50").Value = "'0.18%"
